$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 values
$ws.Range("B2").Value = 49.062128841075321
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 52.205685730939273
$ws.Range("E2").Value = 57.78619717019599

# Row 3 values
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 48.282428813418115
$ws.Range("D3").Value = 46.926664468428058
$ws.Range("E3").Value = 56.091780148283107

# Update selection to match new used region
$ws.Range("B1:E3").Select() | Out-Null
